$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AT")
$c = $ws.Range("I41")
Write-Host "I41 font color:" $c.Font.Color
Write-Host "I41 font underline:" $c.Font.Underline
Write-Host "I41 borders edgebottom:" $c.Borders.Item(9).LineStyle
$c2 = $ws.Range("I42")
Write-Host "I42 borders edgebottom:" $c2.Borders.Item(9).LineStyle
Write-Host "I42 borders edgetop:" $c2.Borders.Item(8).LineStyle
$c3 = $ws.Range("P42")
Write-Host "P42 value:" $c3.Value2
